$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update time log dates from 08 Apr 2019 (43563) to 20 Apr 2019 (43575)
$ws.Range("C7").Value = 43575
$ws.Range("C8").Value = 43575
$ws.Range("C9").Value = 43575
$ws.Range("C10").Value = 43575
$ws.Range("C11").Value = 43575

# Update iteration label
$ws.Range("B3").Value = "BCPR280-Iteration5"

# Move selection to B3
$ws.Range("B3").Select()
